$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 23 (existing rows 23-41 shift down to 26-44)
$ws.Range('A23:A25').EntireRow.Insert()

# New row 23 - Comercializadora del Agro de Limarí, Coquimbo, Tuna Especial (2022-02-17)
$ws.Range('A23').Value() = 2
$ws.Range('B23').Value() = 'Comercializadora del Agro de Limarí'
$ws.Range('C23').Value() = 'Coquimbo'
$ws.Range('D23').Value() = 44609
$ws.Range('E23').Value() = 4
$ws.Range('F23').Value() = 'Fruta'
$ws.Range('G23').Value() = 100107
$ws.Range('H23').Value() = 'Otros'
$ws.Range('I23').Value() = 100107011
$ws.Range('J23').Value() = 'Tuna'
$ws.Range('K23').Value() = 'Sin especificar'
$ws.Range('L23').Value() = 'Especial'
$ws.Range('M23').Value() = 400
$ws.Range('N23').Value() = 14000
$ws.Range('O23').Value() = 15000
$ws.Range('P23').Value() = 14500
$ws.Range('Q23').Value() = '$/caja 18 kilos'
$ws.Range('R23').Value() = 'Provincia de Limarí'
$ws.Range('S23').Value() = 806
$ws.Range('T23').Value() = 18

# New row 24 - Primera (2022-02-17)
$ws.Range('A24').Value() = 2
$ws.Range('B24').Value() = 'Comercializadora del Agro de Limarí'
$ws.Range('C24').Value() = 'Coquimbo'
$ws.Range('D24').Value() = 44609
$ws.Range('E24').Value() = 4
$ws.Range('F24').Value() = 'Fruta'
$ws.Range('G24').Value() = 100107
$ws.Range('H24').Value() = 'Otros'
$ws.Range('I24').Value() = 100107011
$ws.Range('J24').Value() = 'Tuna'
$ws.Range('K24').Value() = 'Sin especificar'
$ws.Range('L24').Value() = 'Primera'
$ws.Range('M24').Value() = 500
$ws.Range('N24').Value() = 11000
$ws.Range('O24').Value() = 12000
$ws.Range('P24').Value() = 11500
$ws.Range('Q24').Value() = '$/caja 18 kilos'
$ws.Range('R24').Value() = 'Provincia de Limarí'
$ws.Range('S24').Value() = 639
$ws.Range('T24').Value() = 18

# New row 25 - Segunda (2022-02-17)
$ws.Range('A25').Value() = 2
$ws.Range('B25').Value() = 'Comercializadora del Agro de Limarí'
$ws.Range('C25').Value() = 'Coquimbo'
$ws.Range('D25').Value() = 44609
$ws.Range('E25').Value() = 4
$ws.Range('F25').Value() = 'Fruta'
$ws.Range('G25').Value() = 100107
$ws.Range('H25').Value() = 'Otros'
$ws.Range('I25').Value() = 100107011
$ws.Range('J25').Value() = 'Tuna'
$ws.Range('K25').Value() = 'Sin especificar'
$ws.Range('L25').Value() = 'Segunda'
$ws.Range('M25').Value() = 300
$ws.Range('N25').Value() = 8000
$ws.Range('O25').Value() = 9000
$ws.Range('P25').Value() = 8500
$ws.Range('Q25').Value() = '$/caja 18 kilos'
$ws.Range('R25').Value() = 'Provincia de Limarí'
$ws.Range('S25').Value() = 472
$ws.Range('T25').Value() = 18
